$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13
$ws.Range("B3").Value = 4.5
$ws.Range("C3").Value = 12
